$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# Paragraph 4: "Execute " + proofErr(spellStart) + "DatabaseScript.sql" + proofErr(spellEnd)
#   -> single run "Execute DatabaseScript.sql", no proofErr
$xml4 = '<w:p ' + $wNs + ' w14:paraId="6C3ACE94" w14:textId="79753680" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Execute DatabaseScript.sql</w:t></w:r></w:p>'
$d.Paragraphs(4).Range.InsertXML($xml4)

# Paragraph 7: proofErr(spellStart) + "AutoMapper" + proofErr(spellEnd) + " 12.0.1"
#   -> single run "AutoMapper 12.0.1"
$xml7 = '<w:p ' + $wNs + ' w14:paraId="5C865C3D" w14:textId="724A2B94" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>AutoMapper 12.0.1</w:t></w:r></w:p>'
$d.Paragraphs(7).Range.InsertXML($xml7)

# Paragraph 8: proofErr(spellStart,gramStart) + "AutoMapper.Extensions.Microsoft.DependencyInjection" + proofErr(spellEnd,gramEnd) + " 12.0.1"
#   -> single run "AutoMapper.Extensions.Microsoft.DependencyInjection 12.0.1"
$xml8 = '<w:p ' + $wNs + ' w14:paraId="4199A3BD" w14:textId="52CCE79D" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>AutoMapper.Extensions.Microsoft.DependencyInjection 12.0.1</w:t></w:r></w:p>'
$d.Paragraphs(8).Range.InsertXML($xml8)

# Paragraph 9: proofErr(spellStart) + "Microsoft.EntityFrameworkCore" + proofErr(spellEnd) + " 7.0.5"
#   -> single run "Microsoft.EntityFrameworkCore 7.0.5"
$xml9 = '<w:p ' + $wNs + ' w14:paraId="7F9BC758" w14:textId="423197EC" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Microsoft.EntityFrameworkCore 7.0.5</w:t></w:r></w:p>'
$d.Paragraphs(9).Range.InsertXML($xml9)

# Paragraph 11: proofErr(spellStart) + "AutoMapper" + proofErr(spellEnd) + " 12.0.1"
#   -> single run "AutoMapper 12.0.1"
$xml11 = '<w:p ' + $wNs + ' w14:paraId="0B05B0F1" w14:textId="1F10CB9E" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>AutoMapper 12.0.1</w:t></w:r></w:p>'
$d.Paragraphs(11).Range.InsertXML($xml11)

# Paragraph 12: proofErr(spellStart) + "EPPlus" + proofErr(spellEnd) + " 4.5.1"
#   -> single run "EPPlus 4.5.1"
$xml12 = '<w:p ' + $wNs + ' w14:paraId="1D3ABFF9" w14:textId="50C2B488" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>EPPlus 4.5.1</w:t></w:r></w:p>'
$d.Paragraphs(12).Range.InsertXML($xml12)

# Paragraph 13: proofErr(spellStart) + "Microsoft.EntityFrameworkCore" + proofErr(spellEnd) + " 7.0.5"
#   -> single run "Microsoft.EntityFrameworkCore 7.0.5"
$xml13 = '<w:p ' + $wNs + ' w14:paraId="1598B1EF" w14:textId="77777777" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Microsoft.EntityFrameworkCore 7.0.5</w:t></w:r></w:p>'
$d.Paragraphs(13).Range.InsertXML($xml13)

# Paragraph 14: proofErr(spellStart,gramStart) + "Microsoft.EntityFrameworkCore" + ".SqlServer" + proofErr(spellEnd,gramEnd) + " 7.0.5"
#   -> single run "Microsoft.EntityFrameworkCore.SqlServer 7.0.5"
$xml14 = '<w:p ' + $wNs + ' w14:paraId="1DEEEC6D" w14:textId="5EA6C390" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Microsoft.EntityFrameworkCore.SqlServer 7.0.5</w:t></w:r></w:p>'
$d.Paragraphs(14).Range.InsertXML($xml14)

# Paragraph 15: proofErr(spellStart,gramStart) + "Microsoft.Extensions.Configuration" + proofErr(spellEnd,gramEnd)
#   -> single run "Microsoft.Extensions.Configuration" (text unchanged, proofErr removed)
$xml15 = '<w:p ' + $wNs + ' w14:paraId="02669176" w14:textId="39E5400D" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Microsoft.Extensions.Configuration</w:t></w:r></w:p>'
$d.Paragraphs(15).Range.InsertXML($xml15)

# Paragraph 16: "Presentation \ " + proofErr(spellStart) + "FruitSA_Test_OJ" + proofErr(spellEnd)
#   -> single run "Presentation \ FruitSA_Test_OJ"
$xml16 = '<w:p ' + $wNs + ' w14:paraId="13DAD69D" w14:textId="5BD6B0AC" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Presentation \ FruitSA_Test_OJ</w:t></w:r></w:p>'
$d.Paragraphs(16).Range.InsertXML($xml16)

# Paragraph 17: proofErr(spellStart) + "EPPlus" + proofErr(spellEnd) + " 7.7.0"
#   -> single run "EPPlus 7.7.0", PLUS a brand new paragraph after it:
#      "Ensure that FruitSA_Test_OJ is the StartUp project" (ListParagraph, ilvl=1, numId=2)
$xml17 = '<w:p ' + $wNs + ' w14:paraId="21025F13" w14:textId="3775C1CA" w:rsidR="00C11C77" w:rsidRDefault="00C11C77" w:rsidP="00D96E2D"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>EPPlus 7.7.0</w:t></w:r></w:p><w:p ' + $wNs + '><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Ensure that FruitSA_Test_OJ is the StartUp project</w:t></w:r></w:p>'
$d.Paragraphs(17).Range.InsertXML($xml17)
